$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Remove the leftover "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair)
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text starts with $marker and replace its
# whole XML with $newXml (a <w:p>...</w:p> fragment).
# ---------------------------------------------------------------------------
function Set-ParagraphXml($marker, $newXml) {
    $count = $d.Paragraphs.Count
    $i = 1
    while ($i -le $count) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($marker)) {
            $p.Range.InsertXML($newXml)
            return
        }
        $i = $i + 1
    }
}

# ---------------------------------------------------------------------------
# 2. "Discussion items" paragraph -> wrap "How" in proofErr gramStart/gramEnd
# ---------------------------------------------------------------------------
$p5 = '<w:p ' + $wns + ' w:rsidR="00781D96" w:rsidRPr="00781D96" w:rsidRDefault="007312EA" w:rsidP="00F63139">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Discussion items</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>How</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> should we refactor our project so as to be more comprehensible?</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "Discussion items" $p5

# ---------------------------------------------------------------------------
# 3. "Outcomes and assignments" paragraph -> wrap "Decided" in gramStart/
#    gramEnd and "CharacterView" in spellStart/spellEnd
# ---------------------------------------------------------------------------
$p6 = '<w:p ' + $wns + ' w:rsidR="00943B31" w:rsidRDefault="00F63139" w:rsidP="00F63139">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r w:rsidRPr="008E18F0"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>O</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>utcomes and assignments</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Decided</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to make a class Character and a corresponding </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CharacterView</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, for the model and the view respectively.</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "Outcomes and assignments" $p6

# ---------------------------------------------------------------------------
# 4. "Wrap up" paragraph -> wrap "Realise" in proofErr spellStart/spellEnd
# ---------------------------------------------------------------------------
$p7 = '<w:p ' + $wns + ' w:rsidR="00FF5835" w:rsidRPr="00F63139" w:rsidRDefault="000509B8" w:rsidP="00F63139">' +
        '<w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Wrap up</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">Goals for next meeting: </w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">    </w:t></w:r>' +
        '<w:r w:rsidR="00C07587"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Write some tests.</w:t></w:r>' +
        '<w:r w:rsidR="00C07587"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">    </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Realise</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> MVC better, split code up into classes.</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "Wrap up" $p7
